# Add I0 and IF columns to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the formatting of the other headers
# (bold font, thin border on all sides, centered horizontally, top vertical alignment)
# by copying the existing header cell's format (reuses the same style definition).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:I68 and J2:J68 (rows 2 through 68).
$iValues = @(9,7,5,9,9,8,9,8,7,8,8,9,8,7,7,7,9,8,10,7,7,8,8,9,9,9,9,7,9,7,9,9,8,8,8,8,10,8,8,9,9,9,9,9,8,8,8,9,9,7,8,6,8,7,10,8,9,9,9,8,8,9,7,5,5,5,6)
$jValues = @(9,7,5,9,9,8,9,8,7,8,8,9,8,7,7,7,9,8,10,7,7,8,8,9,9,9,9,8,9,7,9,9,8,8,8,8,10,8,8,9,9,9,10,9,8,8,8,9,9,7,8,6,8,8,10,8,9,9,9,8,8,9,7,5,5,5,6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
